# Auto-generated edit script: append a new COVID-19 data block
# (header row 649 + 35 state/UT rows 650-684, dated 05-10-2020)
# to the COVID19_TIMESERIESDATA worksheet, matching the pattern of the
# repeated 36-row blocks already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 649) - same headers/formatting as the other block
#    headers in the sheet (bold, centered, top-aligned, thin box border).
# ---------------------------------------------------------------------------
$headerRow = 649
$headers = @('States/UT', 'Active Cases', 'Active Cases Since Yesterday', 'Recovered Cases', 'Recovered Cases Since Yesterday', 'Deceased Cases', 'Deceased Cases Since Yesterday', 'Date')

for ($col = 1; $col -le 8; $col++) {
    $cell = $ws.Cells.Item($headerRow, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous (thin box border)
}

# ---------------------------------------------------------------------------
# 2. Data rows (rows 650-684): one row per State/UT for date 05-10-2020.
#    Columns: A=States/UT, B=Active Cases, C=Active Cases Since Yesterday,
#    D=Recovered Cases, E=Recovered Cases Since Yesterday, F=Deceased Cases,
#    G=Deceased Cases Since Yesterday, H=Date.
# ---------------------------------------------------------------------------
$data = @(
    @('Andaman and Nicobar Islands', 182, 9, 3649, 7, 53, 0, '05-10-2020'),
    @('Andhra Pradesh', 54400, -882, 658875, 7084, 5981, 40, '05-10-2020'),
    @('Arunachal Pradesh', 2953, -62, 7577, 189, 18, 0, '05-10-2020'),
    @('Assam', 33324, -609, 152127, 1351, 749, 14, '05-10-2020'),
    @('Bihar', 11795, 198, 175458, 1526, 915, 3, '05-10-2020'),
    @('Chandigarh', 1673, -119, 10598, 202, 174, 2, '05-10-2020'),
    @('Chhattisgarh', 28548, -744, 93731, 2654, 1045, 14, '05-10-2020'),
    @('Dadra and Nagar Haveli and Daman and Diu', 105, 4, 2980, 13, 2, 0, '05-10-2020'),
    @('Delhi', 24753, -481, 260350, 3126, 5510, 38, '05-10-2020'),
    @('Goa', 4839, -84, 30033, 506, 456, 6, '05-10-2020'),
    @('Gujarat', 16809, 47, 122233, 1246, 3496, 9, '05-10-2020'),
    @('Haryana', 12067, -801, 120341, 2083, 1470, 20, '05-10-2020'),
    @('Himachal Pradesh', 3273, -19, 12361, 168, 217, 7, '05-10-2020'),
    @('Jammu and Kashmir', 15460, -186, 62404, 1053, 1242, 11, '05-10-2020'),
    @('Jharkhand', 10936, -3, 75531, 927, 743, 9, '05-10-2020'),
    @('Karnataka', 115593, 2791, 515782, 7287, 9286, 67, '05-10-2020'),
    @('Kerala', 84579, 3679, 144471, 4851, 836, 23, '05-10-2020'),
    @('Ladakh', 1106, 5, 3354, 39, 61, 0, '05-10-2020'),
    @('Madhya Pradesh', 19372, -435, 113832, 2120, 2434, 35, '05-10-2020'),
    @('Maharashtra', 255722, -2826, 1149603, 15048, 38084, 326, '05-10-2020'),
    @('Manipur', 2576, 79, 9205, 124, 74, 3, '05-10-2020'),
    @('Meghalaya', 2209, 126, 4393, 74, 54, 0, '05-10-2020'),
    @('Mizoram', 313, -35, 1807, 35, 0, 0, '05-10-2020'),
    @('Nagaland', 1226, 71, 5309, 52, 17, 0, '05-10-2020'),
    @('Odisha', 29504, -797, 202302, 4108, 907, 15, '05-10-2020'),
    @('Puducherry', 4787, -87, 23763, 419, 539, 5, '05-10-2020'),
    @('Punjab', 13577, -712, 100977, 1509, 3603, 41, '05-10-2020'),
    @('Rajasthan', 21154, 79, 121331, 2090, 1545, 15, '05-10-2020'),
    @('Sikkim', 649, 9, 2480, 32, 45, 2, '05-10-2020'),
    @('Tamil Nadu', 46120, -135, 564092, 5558, 9784, 66, '05-10-2020'),
    @('Telengana', 27052, -849, 172388, 2176, 1171, 8, '05-10-2020'),
    @('Tripura', 4858, -313, 21876, 466, 299, 6, '05-10-2020'),
    @('Uttarakhand', 9089, 1013, 41740, 402, 652, 4, '05-10-2020'),
    @('Uttar Pradesh', 46385, -1438, 362052, 5226, 6029, 52, '05-10-2020'),
    @('West Bengal', 27439, 309, 237698, 2986, 5194, 62, '05-10-2020'),
)

$startRow = 650

# Column H holds a dd-mm-yyyy formatted date stored as literal text (matching
# every other "Date" column entry in this sheet) - format the range as Text
# first so Excel does not reinterpret the string as a date serial number.
$ws.Range("H" + $startRow + ":H" + ($startRow + $data.Count - 1)).NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
    $ws.Cells.Item($r, 6).Value = $rowVals[5]
    $ws.Cells.Item($r, 7).Value = $rowVals[6]
    $ws.Cells.Item($r, 8).Value = $rowVals[7]
}
